# The workbook already holds a MAG rejection-prediction table. Previously the
# header row read: max | prediction | rejection-f (columns C, D, E) and the
# "max" column (C) duplicated the numeric "max-score" from column B, while
# "rejection-f" (E) duplicated the predicted taxon string from D.
#
# The update reorders the headers to: prediction | rejection-f | max, and
# updates the data rows so that column C now holds the predicted taxon
# string (same text as D), and column E now holds the numeric rejection
# flag (1) instead of duplicating the taxon string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = "s__QAMX01 sp003149835"   # column C
    $ws.Cells.Item($r, 5).Value = 1                          # column E
}
